$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (single decimal point) need to be
# forced to text so Excel does not silently convert them to numbers.
$textCells = @("D5", "D6", "D7", "D10", "D12", "D14", "D16", "D18", "D19", "D21", "D23", "D24", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D41", "D42", "D43", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "188.48"
$ws.Range("D6").Value = "557.47"
$ws.Range("D7").Value = "1.00"
$ws.Range("D10").Value = "0.186"
$ws.Range("D12").Value = "47.61"
$ws.Range("D14").Value = "8.67"
$ws.Range("D16").Value = "626.57"
$ws.Range("D18").Value = "18.08"
$ws.Range("D19").Value = "0.118"
$ws.Range("D21").Value = "10.84"
$ws.Range("D23").Value = "18.14"
$ws.Range("D24").Value = "102.98"
$ws.Range("D29").Value = "9.63"
$ws.Range("D30").Value = "8.70"
$ws.Range("D31").Value = "30.37"
$ws.Range("D32").Value = "4.02"
$ws.Range("D33").Value = "6.40"
$ws.Range("D34").Value = "560.01"
$ws.Range("D36").Value = "0.106"
$ws.Range("D38").Value = "58.01"
$ws.Range("D39").Value = "0.999"
$ws.Range("D41").Value = "34.16"
$ws.Range("D42").Value = "3.32"
$ws.Range("D43").Value = "2.73"
$ws.Range("D49").Value = "0.130"
$ws.Range("D50").Value = "2.60"
$ws.Range("D51").Value = "1.00"

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining text/percentage/link cells - plain string assignment is safe.
$ws.Range("D2").Value = "66.290.26"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "3.314.71"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").Value = "3.306.09"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "3.845.83"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "66.317.02"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").Value = "3.320.33"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("E23").Value = "  +6.28%  "
$ws.Range("E24").Value = "  +6.20%  "
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.850.80"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "0.0₃0733"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E41").Value = "  +3.70%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  -4.68%  "
$ws.Range("E46").Value = "  -15.55%  "
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("E50").Value = "  -3.55%  "
$ws.Range("E51").Value = "  +0.04%  "
